$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (s="1") from an existing header cell (H1) onto I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-26
$values = @(
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(6, 6),
    @(6, 7),
    @(4, 5),
    @(5, 6),
    @(3, 3),
    @(11, 11),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(1, 1),
    @(6, 6),
    @(2, 2),
    @(4, 4),
    @(3, 3),
    @(1, 1)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
